$wb = $excel.ActiveWorkbook

# --- 1. Update the JMP (Job Market Paper) abstract on "Working Papers" ---
$wsWorking = $wb.Worksheets.Item("Working Papers")
$wsWorking.Range("C2").Value = "Economists typically view firms as risk neutral. Yet many enterprises, especially in low and middle-income economies, are small and owner-operated, making household consumption sensitive to business risk. As a result, owners' risk preferences may influence firm decisions. This paper demonstrates that small retailers in Kenya are risk averse, leading them to under-adopt a new product when they face uncertain demand. I model risk averse firms who learn about demand through stocking decisions, then test the model's predictions using two field experiments. The first establishes that risk aversion affects the stocking decisions of enterprises. I test for risk aversion by offering treated firms an insurance contract that lowers expected profits from a new product while reducing the risk of losses. This leads to a 50\% increase in adoption, rejecting risk neutrality. The second experiment shows that \textit{temporarily} reducing inventory risk leads firms to \textit{permanently} stock a profitable new product because they overcome demand uncertainty through learning. These results show that risk aversion in firms can impede product diffusion, potentially limiting growth."

# --- 2. Add the new "Long-run Effects of Unconditional Cash Transfers" row to "Work in Progress" ---
$wsWip = $wb.Worksheets.Item("Work in Progress")
$wsWip.Range("A4").Value = "The Long-run Effects of Unconditional Cash Transfers: Evidence from the Kenya General Equilibrium Study"
$wsWip.Range("B4").Value = "with Dennis Egger, Edward Miguel and Michael Walker"
$wsWip.Range("C4").Value = "Recent studies document positive short-run effects of unconditional cash transfers (UCTs) on cash recipients and spillover effects on non-recipients. But modest sample sizes and challenges with tracking households over time have limited research on the long-run effects of UCTs. We study the long-run effects of the Kenya General Equilibrium Study (KGES) on recipient households and the local economy using census and survey data collected up to ten years post-transfers. An “endline 2” completed 4-7 years after the transfers collected census data from each household and enterprise in the study area, and a representative survey of more than 10,000 households and firms obtained more detailed consumption and production measures for a subset, including from those that migrated out of the study area. We collected this same information in an “endline 3” completed 7-10 years after the experimental start. Preliminary results from endline 2 show persistent consumption gains among recipients and expansion of non-farm enterprise revenue. Ongoing analysis examines whether these effects persisted to endline 3 and aims to estimate long-run transfer multipliers. "

# --- 3. Update selections / active sheet to match final view state ---
$wsWorking.Range("C3").Select() | Out-Null

$wsWip.Activate()
$wsWip.Range("C5").Select() | Out-Null
